$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Plan1) cosmetic changes -------------------------------------
# Widen column D to fit the new/longer product names
$ws1.Columns.Item(4).ColumnWidth = 27.140625

# Update the view: drop the frozen/top-left scroll position and move the
# active selection
$ws1.Range("G20").Select()

# --- Sheet2 (Plan2) gains a distinct BIN list -----------------------------
# Header row, copied from Plan1 (keeps formatting + shared strings intact)
$ws1.Range("A1:D1").Copy($ws2.Range("A1:D1"))

# Source rows on Plan1 holding the first occurrence of each distinct
# ORG/PRODUTO combination (columns A-D only, no BIN/CHIP columns)
$sourceRows = @(2, 3, 12, 26, 33, 34, 37, 39, 41, 42, 44, 45, 46, 60, 64, 65, 66, 67, 68, 69, 70, 71, 73, 82, 83, 84)

$destRow = 2
foreach ($srcRow in $sourceRows) {
    $srcRange = $ws1.Range("A" + $srcRow + ":D" + $srcRow)
    $dstRange = $ws2.Range("A" + $destRow + ":D" + $destRow)
    $srcRange.Copy($dstRange)
    $destRow++
}

# Match column widths with Plan1
$ws2.Columns.Item(3).ColumnWidth = 14.5703125
$ws2.Columns.Item(4).ColumnWidth = 27.140625

# Set the new active selection on Plan2
$ws2.Range("F14").Select()
